# Camera Bazaar edit profile finished
#
# The "Edit profile" requirement paragraph ("Each user can change his
# email, phone or password.") is marked as done: a bold "[done]" run is
# appended at the end of the paragraph (after the trailing space run,
# still inside the same <w:p>), wrapped in a _GoBack bookmark. Word only
# ever keeps a single "_GoBack" bookmark in a document, so adding the new
# one here automatically relocates it away from its old position later in
# the document (after the "Picture 10" drawing), matching the target diff.

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute(
    "Each user can change his email, phone or password.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the 'edit profile' requirement paragraph"
}

$para = $rng.Paragraphs(1)
$paraRange = $para.Range

# Paragraph Range.End sits just past the paragraph mark, so back up one
# character to get a collapsed insertion point right before it (i.e.
# after the existing trailing-space run, still inside this <w:p>).
$insertAt = $paraRange.End - 1

$insertPoint = $d.Range($insertAt, $insertAt)
$insertPoint.InsertAfter("[done]")

$doneRange = $d.Range($insertAt, $insertAt + 6)
$doneRange.Bold = 1

$d.Bookmarks.Add("_GoBack", $doneRange) | Out-Null
